$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": insert new row 253 ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows("253:253").Insert()
$ws1.Range("A253").Value = "OFICINA-CATAECSA"
$ws1.Range("B253").Value = "SOLIS OCAMPO DIMAS ABDON"
$ws1.Range("C253:R253").Value = 0

# Update the "x de N" summary row (now shifted from 280 to 281) to reflect the new row count
$ws1.Range("C281").Value = "0 de 279"
$ws1.Range("D281").Value = "0 de 279"
$ws1.Range("E281").Value = "0 de 279"
$ws1.Range("F281").Value = "0 de 279"
$ws1.Range("G281").Value = "0 de 279"
$ws1.Range("H281").Value = "0 de 279"
$ws1.Range("I281").Value = "0 de 279"
$ws1.Range("J281").Value = "0 de 279"
$ws1.Range("K281").Value = "1 de 279"
$ws1.Range("L281").Value = "0 de 279"
$ws1.Range("M281").Value = "1 de 279"
$ws1.Range("N281").Value = "0 de 279"
$ws1.Range("O281").Value = "0 de 279"
$ws1.Range("P281").Value = "0 de 279"
$ws1.Range("Q281").Value = "1 de 279"
$ws1.Range("R281").Value = "0 de 279"

# --- Sheet "VENTA MENSUAL": insert new row 253 ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows("253:253").Insert()
$ws2.Range("A253").Value = "OFICINA-CATAECSA"
$ws2.Range("B253").Value = "SOLIS OCAMPO DIMAS ABDON"
$ws2.Range("C253:G253").Value = 0
